$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (Id=5): Name/Cname/Icon/Color take on what used to be row 10's values
$ws.Range("B8").Value = "light"
$ws.Range("C8").Value = "遗迹"
$ws.Range("G8").Value = "Gold"

# Row 9 (Id=6): Name/Cname/Icon/Color take on what used to be row 11's values
$ws.Range("B9").Value = "dark"
$ws.Range("C9").Value = "沼泽"
$ws.Range("G9").Value = "Brown"

# Row 10 (Id=7): Name/Cname/Icon/Color take on what used to be row 8's values; Type becomes 0
$ws.Range("B10").Value = "snow"
$ws.Range("C10").Value = "雪地"
$ws.Range("D10").Value = 0
$ws.Range("G10").Value = "White"

# Row 11 (Id=8): Name/Cname/Icon/Color take on what used to be row 9's values; Type becomes 0
$ws.Range("B11").Value = "hill"
$ws.Range("C11").Value = "山地"
$ws.Range("D11").Value = 0
$ws.Range("G11").Value = "DarkSlateGray"

# F column (Icon) cells are formatted as Text ("@"); writing .Value directly would
# coerce the number to a text shared-string. Clear formats, write the numeric
# value, then restore the original Text format by copying it from a sibling
# cell that already carries it (avoids minting a new, unused style entry).
function Set-NumericTextCell($cell, $value, $formatSource) {
    $cell.ClearFormats()
    $cell.Value = $value
    $formatSource.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

Set-NumericTextCell $ws.Range("F8") 7 $ws.Range("F4")
Set-NumericTextCell $ws.Range("F9") 8 $ws.Range("F4")
Set-NumericTextCell $ws.Range("F10") 5 $ws.Range("F4")
Set-NumericTextCell $ws.Range("F11") 6 $ws.Range("F4")

# Selection moved from B13 to B5
$ws.Range("B5").Select()
